$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K8").Value = 0
$ws.Range("M10").Value = 10631.67
$ws.Range("M11").Value = 790500.47
$ws.Range("O11").Value = 333891.59
$ws.Range("M12").Value = 135737.98
$ws.Range("O12").Value = 45523.6
$ws.Range("O26").Value = 43190
